# Update database and change read_price algorithm
# - shifts each yearly column one period to the left (drop oldest, keep
#   the four most-recent periods) and appends a new period's figures
# - updates the period-label / publish-date header rows to match

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "12 ماهه منتهی به ..." period headers (D8:H8) ---
$ws.Range("D8").Value = "12 ماهه منتهی به 1397/12"
$ws.Range("E8").Value = "12 ماهه منتهی به 1398/12"
$ws.Range("F8").Value = "12 ماهه منتهی به 1399/12"
$ws.Range("G8").Value = "12 ماهه منتهی به 1400/12"
$ws.Range("H8").Value = "12 ماهه منتهی به 1401/12"

# --- Row 9: "تاریخ انتشار" publish-date headers (D9:H9) ---
$ws.Range("D9").Value = "1399-04-08 (8)"
$ws.Range("E9").Value = "1400-04-20 (8)"
$ws.Range("F9").Value = "1401-04-26 (9)"
$ws.Range("G9").Value = "1402-02-30 (8)"
$ws.Range("H9").Value = "1402-02-30"

# --- Financial data rows: shift D<-E<-F<-G<-H, fill H with the new figure ---
$cols = @("D", "E", "F", "G", "H")

$rowsData = @{
    11 = @(29185, 34329, 22396, 40389, 41717)
    12 = @(-24651, -28082, -18854, -34875, -37289)
    13 = @(4535, 6247, 3542, 5514, 4428)
    14 = @(-659, -656, -442, -674, -753)
    16 = @(690, -24, -3, -37, 0)
    17 = @(4566, 5567, 3098, 4802, 3675)
    18 = @(-642, -840, -96, -125, -99)
    19 = @(1784, 193, 351, 77, 688)
    20 = @(5709, 4920, 3353, 4754, 4263)
    21 = @(-991, -979, -646, -691, -636)
    22 = @(4718, 3941, 2707, 4063, 3628)
    24 = @(4718, 3941, 2707, 4063, 3628)
    26 = @(1928, 1520, 862, 738, 552)
}

foreach ($row in $rowsData.Keys) {
    $values = $rowsData[$row]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $addr = $cols[$i] + $row
        $ws.Range($addr).Value = $values[$i]
    }
}

Write-Host "Applied database update + read_price shift."
